$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Page1")

$ws.Range("A2").Value = "Contract python test developer"
$ws.Range("B2").Value = "https://uk.indeed.com/rc/clk?jk=705320533dc92838&fccid=8f35589a37e69470&vjs=3"

$ws.Range("A3").Value = "Software Engineer in Test (SDET)"
$ws.Range("B3").Value = "https://uk.indeed.com/rc/clk?jk=fd9e8860b3959cc4&fccid=c659788ec6cc356e&vjs=3"
